$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 1.296878333333333
$ws.Range("H2").Value = 3.890635
$ws.Range("I2").Value = 0.01774073260139904
$ws.Range("J2").Value = 0.02506266560199287
$ws.Range("M2").Value = 0.169654
$ws.Range("N2").Value = 0.508962
$ws.Range("O2").Value = 0.006094264463659866
$ws.Range("P2").Value = 0.006534681579452628
$ws.Range("Q2").Value = 0.2200205967633333
$ws.Range("R2").Value = 1.98018537087
$ws.Range("S2").Value = 0.0001081167162519982
$ws.Range("T2").Value = 0.0001637765392413238
$ws.Range("G3").Value = 1.296878333333333
$ws.Range("H3").Value = 3.890635
$ws.Range("I3").Value = 0.01774073260139904
$ws.Range("J3").Value = 0.02506266560199287
$ws.Range("O3").Value = 0.7900017288527916
$ws.Range("P3").Value = 0.8470931604713817
$ws.Range("Q3").Value = 28.52135033895
$ws.Range("R3").Value = 256.69215305055
$ws.Range("S3").Value = 0.01401520942622032
$ws.Range("T3").Value = 0.02123041261462952
$ws.Range("G4").Value = 1.296878333333333
$ws.Range("H4").Value = 3.890635
$ws.Range("I4").Value = 0.01774073260139904
$ws.Range("J4").Value = 0.02506266560199287
$ws.Range("M4").Value = 0.04769766666666667
$ws.Range("N4").Value = 0.143093
$ws.Range("O4").Value = 0.001713382501834088
$ws.Range("P4").Value = 0.001837204332049573
$ws.Range("Q4").Value = 0.06185807045055555
$ws.Range("R4").Value = 0.556722634055
$ws.Range("S4").Value = 0.00003039666080895465
$ws.Range("T4").Value = 0.00004604523781669112
$ws.Range("G5").Value = 1.296878333333333
$ws.Range("H5").Value = 3.890635
$ws.Range("I5").Value = 0.01774073260139904
$ws.Range("J5").Value = 0.02506266560199287
$ws.Range("M5").Value = 5.6286445
$ws.Range("N5").Value = 11.257289
$ws.Range("O5").Value = 0.2021906241817143
$ws.Range("P5").Value = 0.1445349536171162
$ws.Range("Q5").Value = 7.299667098085833
$ws.Range("R5").Value = 43.798002588515
$ws.Range("S5").Value = 0.003587009798117759
$ws.Range("T5").Value = 0.003622431210305333
$ws.Range("I6").Value = 0.09943605305674341
$ws.Range("J6").Value = 0.1404751766759988
$ws.Range("M6").Value = 0.169654
$ws.Range("N6").Value = 0.508962
$ws.Range("O6").Value = 0.006094264463659866
$ws.Range("P6").Value = 0.006534681579452628
$ws.Range("Q6").Value = 1.233206103992
$ws.Range("R6").Value = 11.098854935928
$ws.Range("S6").Value = 0.0006059896045503084
$ws.Range("T6").Value = 0.0009179605493950027
$ws.Range("I7").Value = 0.09943605305674341
$ws.Range("J7").Value = 0.1404751766759988
$ws.Range("O7").Value = 0.7900017288527916
$ws.Range("P7").Value = 0.8470931604713817
$ws.Range("S7").Value = 0.07855465382512521
$ws.Range("T7").Value = 0.1189955613782475
$ws.Range("I8").Value = 0.09943605305674341
$ws.Range("J8").Value = 0.1404751766759988
$ws.Range("M8").Value = 0.04769766666666667
$ws.Range("N8").Value = 0.143093
$ws.Range("O8").Value = 0.001713382501834088
$ws.Range("P8").Value = 0.001837204332049573
$ws.Range("Q8").Value = 0.3467118587213333
$ws.Range("R8").Value = 3.120406728492
$ws.Range("S8").Value = 0.0001703719933588702
$ws.Range("T8").Value = 0.0002580816031345741
$ws.Range("I9").Value = 0.09943605305674341
$ws.Range("J9").Value = 0.1404751766759988
$ws.Range("M9").Value = 5.6286445
$ws.Range("N9").Value = 11.257289
$ws.Range("O9").Value = 0.2021906241817143
$ws.Range("P9").Value = 0.1445349536171162
$ws.Range("Q9").Value = 40.91432418098599
$ws.Range("R9").Value = 245.485945085916
$ws.Range("S9").Value = 0.02010503763370901
$ws.Range("T9").Value = 0.02030357314522169
$ws.Range("G10").Value = 0.09795233333333332
$ws.Range("H10").Value = 0.293857
$ws.Range("I10").Value = 0.001339945397100812
$ws.Range("J10").Value = 0.001892965987764162
$ws.Range("M10").Value = 0.169654
$ws.Range("N10").Value = 0.508962
$ws.Range("O10").Value = 0.006094264463659866
$ws.Range("P10").Value = 0.006534681579452628
$ws.Range("Q10").Value = 0.01661800515933333
$ws.Range("R10").Value = 0.149562046434
$ws.Range("S10").Value = 0.000008165981616796084
$ws.Range("T10").Value = 0.00001236992997077282
$ws.Range("G11").Value = 0.09795233333333332
$ws.Range("H11").Value = 0.293857
$ws.Range("I11").Value = 0.001339945397100812
$ws.Range("J11").Value = 0.001892965987764162
$ws.Range("O11").Value = 0.7900017288527916
$ws.Range("P11").Value = 0.8470931604713817
$ws.Range("Q11").Value = 2.15419807989
$ws.Range("R11").Value = 19.38778271901
$ws.Range("S11").Value = 0.001058559180277981
$ws.Range("T11").Value = 0.001603518541239975
$ws.Range("G12").Value = 0.09795233333333332
$ws.Range("H12").Value = 0.293857
$ws.Range("I12").Value = 0.001339945397100812
$ws.Range("J12").Value = 0.001892965987764162
$ws.Range("M12").Value = 0.04769766666666667
$ws.Range("N12").Value = 0.143093
$ws.Range("O12").Value = 0.001713382501834088
$ws.Range("P12").Value = 0.001837204332049573
$ws.Range("Q12").Value = 0.004672097744555555
$ws.Range("R12").Value = 0.042048879701
$ws.Range("S12").Value = 0.000002295838996805659
$ws.Range("T12").Value = 0.000003477765313142817
$ws.Range("G13").Value = 0.09795233333333332
$ws.Range("H13").Value = 0.293857
$ws.Range("I13").Value = 0.001339945397100812
$ws.Range("J13").Value = 0.001892965987764162
$ws.Range("M13").Value = 5.6286445
$ws.Range("N13").Value = 11.257289
$ws.Range("O13").Value = 0.2021906241817143
$ws.Range("P13").Value = 0.1445349536171162
$ws.Range("Q13").Value = 0.5513388622788332
$ws.Range("R13").Value = 3.308033173673
$ws.Range("S13").Value = 0.0002709243962092281
$ws.Range("T13").Value = 0.0002735997512402716
$ws.Range("G14").Value = 64.06892400000001
$ws.Range("H14").Value = 128.137848
$ws.Range("I14").Value = 0.8764350668284411
$ws.Range("J14").Value = 0.8254375019458241
$ws.Range("M14").Value = 0.169654
$ws.Range("N14").Value = 0.508962
$ws.Range("O14").Value = 0.006094264463659866
$ws.Range("P14").Value = 0.006534681579452628
$ws.Range("Q14").Value = 10.869549232296
$ws.Range("R14").Value = 65.21729539377601
$ws.Range("S14").Value = 0.005341227082477929
$ws.Range("T14").Value = 0.005393971238954769
$ws.Range("G15").Value = 64.06892400000001
$ws.Range("H15").Value = 128.137848
$ws.Range("I15").Value = 0.8764350668284411
$ws.Range("J15").Value = 0.8254375019458241
$ws.Range("O15").Value = 0.7900017288527916
$ws.Range("P15").Value = 0.8470931604713817
$ws.Range("Q15").Value = 1409.02363797444
$ws.Range("R15").Value = 8454.141827846641
$ws.Range("S15").Value = 0.6923852180216804
$ws.Range("T15").Value = 0.6992224622948905
$ws.Range("G16").Value = 64.06892400000001
$ws.Range("H16").Value = 128.137848
$ws.Range("I16").Value = 0.8764350668284411
$ws.Range("J16").Value = 0.8254375019458241
$ws.Range("M16").Value = 0.04769766666666667
$ws.Range("N16").Value = 0.143093
$ws.Range("O16").Value = 0.001713382501834088
$ws.Range("P16").Value = 0.001837204332049573
$ws.Range("Q16").Value = 3.055938180644
$ws.Range("R16").Value = 18.335629083864
$ws.Range("S16").Value = 0.001501668507497641
$ws.Range("T16").Value = 0.001516497354411046
$ws.Range("G17").Value = 64.06892400000001
$ws.Range("H17").Value = 128.137848
$ws.Range("I17").Value = 0.8764350668284411
$ws.Range("J17").Value = 0.8254375019458241
$ws.Range("M17").Value = 5.6286445
$ws.Range("N17").Value = 11.257289
$ws.Range("O17").Value = 0.2021906241817143
$ws.Range("P17").Value = 0.1445349536171162
$ws.Range("Q17").Value = 360.6211966935181
$ws.Range("R17").Value = 1442.484786774072
$ws.Range("S17").Value = 0.177206953216785
$ws.Range("T17").Value = 0.1193045710575679
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.3690323333333334
$ws.Range("H18").Value = 1.107097
$ws.Range("I18").Value = 0.005048202116315478
$ws.Range("J18").Value = 0.007131689788420014
$ws.Range("M18").Value = 0.169654
$ws.Range("N18").Value = 0.508962
$ws.Range("O18").Value = 0.006094264463659866
$ws.Range("P18").Value = 0.006534681579452628
$ws.Range("Q18").Value = 0.06260781147933334
$ws.Range("R18").Value = 0.5634703033140001
$ws.Range("S18").Value = 0.00003076507876283395
$ws.Range("T18").Value = 0.00004660332189075868
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.3690323333333334
$ws.Range("H19").Value = 1.107097
$ws.Range("I19").Value = 0.005048202116315478
$ws.Range("J19").Value = 0.007131689788420014
$ws.Range("O19").Value = 0.7900017288527916
$ws.Range("P19").Value = 0.8470931604713817
$ws.Range("Q19").Value = 8.11587347469
$ws.Range("R19").Value = 73.04286127220999
$ws.Range("S19").Value = 0.003988088399487549
$ws.Range("T19").Value = 0.006041205642374189
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.3690323333333334
$ws.Range("H20").Value = 1.107097
$ws.Range("I20").Value = 0.005048202116315478
$ws.Range("J20").Value = 0.007131689788420014
$ws.Range("M20").Value = 0.04769766666666667
$ws.Range("N20").Value = 0.143093
$ws.Range("O20").Value = 0.001713382501834088
$ws.Range("P20").Value = 0.001837204332049573
$ws.Range("Q20").Value = 0.01760198122455556
$ws.Range("R20").Value = 0.158417831021
$ws.Range("S20").Value = 0.000008649501171816752
$ws.Range("T20").Value = 0.00001310237137411895
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.3690323333333334
$ws.Range("H21").Value = 1.107097
$ws.Range("I21").Value = 0.005048202116315478
$ws.Range("J21").Value = 0.007131689788420014
$ws.Range("M21").Value = 5.6286445
$ws.Range("N21").Value = 11.257289
$ws.Range("O21").Value = 0.2021906241817143
$ws.Range("P21").Value = 0.1445349536171162
$ws.Range("Q21").Value = 2.077151813338833
$ws.Range("R21").Value = 12.462910880033
$ws.Range("S21").Value = 0.001020699136893278
$ws.Range("T21").Value = 0.001030778452780948
